$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.094.69'
$ws.Cells.Item(2, 5).Value = '  -0.99%  '
$ws.Cells.Item(3, 4).Value = '1.822.42'
$ws.Cells.Item(3, 5).Value = '  -1.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.011'
$ws.Cells.Item(4, 5).Value = '  -0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '311.81'
$ws.Cells.Item(5, 5).Value = '  -1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.009'
$ws.Cells.Item(6, 5).Value = '  -0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4629'
$ws.Cells.Item(7, 5).Value = '  -2.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3639'
$ws.Cells.Item(8, 5).Value = '  -1.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07300'
$ws.Cells.Item(9, 5).Value = '  -2.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.8698'
$ws.Cells.Item(10, 5).Value = '  -1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '20.10'
$ws.Cells.Item(11, 5).Value = '  -2.01%  '
$ws.Cells.Item(12, 4).Value = '1.867.76'
$ws.Cells.Item(12, 5).Value = '  +1.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.07597'
$ws.Cells.Item(13, 5).Value = '  +3.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.344'
$ws.Cells.Item(14, 5).Value = '  -2.59%  '
$ws.Cells.Item(15, 5).Value = '  -1.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '6.481'
$ws.Cells.Item(17, 5).Value = '  -0.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000008639'
$ws.Cells.Item(18, 5).Value = '  -2.27%  '
$ws.Cells.Item(19, 5).Value = '  -0.27%  '
$ws.Cells.Item(20, 4).Value = '27.378.78'
$ws.Cells.Item(20, 5).Value = '  -0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '14.46'
$ws.Cells.Item(21, 5).Value = '  -2.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.200'
$ws.Cells.Item(22, 5).Value = '  -2.62%  '
$ws.Cells.Item(23, 5).Value = '  -1.64%  '
$ws.Cells.Item(24, 4).Value = '2.092.25'
$ws.Cells.Item(24, 5).Value = '  +0.86%  '
$ws.Cells.Item(25, 2).Value = 'Toncoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '1.876'
$ws.Cells.Item(25, 5).Value = '  -1.69%  '
$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '151.84'
$ws.Cells.Item(26, 5).Value = '  -0.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '18.24'
$ws.Cells.Item(27, 5).Value = '  -2.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.092'
$ws.Cells.Item(28, 5).Value = '  -4.46%  '
$ws.Cells.Item(29, 2).Value = 'BitcoinCash'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '116.20'
$ws.Cells.Item(29, 5).Value = '  -1.62%  '
$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '5.090'
$ws.Cells.Item(30, 5).Value = '  -3.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.08910'
$ws.Cells.Item(31, 5).Value = '  -0.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '2.950'
$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.7326'
$ws.Cells.Item(33, 5).Value = '  -3.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.454'
$ws.Cells.Item(34, 5).Value = '  -2.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.137'
$ws.Cells.Item(35, 5).Value = '  -3.58%  '
$ws.Cells.Item(36, 5).Value = '  -0.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.475'
$ws.Cells.Item(37, 5).Value = '  +2.21%  '
$ws.Cells.Item(38, 5).Value = '  -3.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.05250'
$ws.Cells.Item(39, 5).Value = '  -2.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.01915'
$ws.Cells.Item(40, 5).Value = '  -2.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.928'
$ws.Cells.Item(41, 5).Value = '  -2.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '7.153'
$ws.Cells.Item(42, 5).Value = '  -2.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.5207'
$ws.Cells.Item(43, 5).Value = '  -2.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.1630'
$ws.Cells.Item(44, 5).Value = '  -2.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '8.261'
$ws.Cells.Item(45, 5).Value = '  -3.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.4878'
$ws.Cells.Item(46, 5).Value = '  -1.91%  '
$ws.Cells.Item(47, 5).Value = '  -0.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '103.75'
$ws.Cells.Item(48, 5).Value = '  -0.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '10.12'
$ws.Cells.Item(49, 5).Value = '  -4.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.630'
$ws.Cells.Item(50, 5).Value = '  -3.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06257'
$ws.Cells.Item(51, 5).Value = '  -1.12%  '
